# Juno: check in to OLPRODLOC.
# Localizes the "Charger sales report" workbook from English to Italian:
#   - renames the worksheet
#   - translates the header row and the Year-Quarter labels in column A

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab + workbook.xml <sheet name=.../>)
$ws.Name = "Report delle vendite"

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Anno-Trimestre"
$ws.Range("B1").Value = "Midwest"
$ws.Range("C1").Value = "Mountain"
$ws.Range("D1").Value = "Northeast"
$ws.Range("E1").Value = "South"
$ws.Range("F1").Value = "sud-orientale"
$ws.Range("G1").Value = "West"

# --- Year-Quarter labels (column A, rows 2-9) -----------------------------
$ws.Range("A2").Value = "2022-T1"
$ws.Range("A3").Value = "2022-T2"
$ws.Range("A4").Value = "2022-T3"
$ws.Range("A5").Value = "2022-T4"
$ws.Range("A6").Value = "2023-T1"
$ws.Range("A7").Value = "2023-T2"
$ws.Range("A8").Value = "2023-T3"
$ws.Range("A9").Value = "2023-T4"
